$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated emission summary results for all basins (PERMIAN output)
$ws.Range("B2").Value = 0.014757983694818918
$ws.Range("C2").Value = 0.05239820010929571

$ws.Range("C3").Value = 0.00043283719946499947

$ws.Range("B4").Value = 0.0011743833068070092
$ws.Range("C4").Value = 0.012054150216426602

$ws.Range("B5").Value = 0.01097501104862853
$ws.Range("C5").Value = 0.038336308322254174

$ws.Range("B6").Value = 0.0046168443208594824

$ws.Range("B7").Value = 0.0021441427380900131
$ws.Range("C7").Value = 0.021593359844230203

$ws.Range("B8").Value = 0.033099041671384652
$ws.Range("C8").Value = 0.16387991142184008

$ws.Range("B9").Value = 0.00088491818994149976

$ws.Range("B10").Value = 0.0000014885263925000002

$ws.Range("B11").Value = 0.0041494499440150297
$ws.Range("C11").Value = 0.020629619761904621

$ws.Range("B12").Value = 0.030761016387854435
$ws.Range("C12").Value = 0.12192541730105945

$ws.Range("B13").Value = 0.0000144213982

$ws.Range("B16").Value = 0.013930385076390757

$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 0.04975552809545003

$ws.Range("B18").Value = 0.10318130996133078
$ws.Range("C18").Value = 0.16611140517099435
